# Refresh the "cryptos" price/volume table with the latest scrape.
# Price cells (column D) that look like a plain decimal number (e.g. "9.50")
# are written with a leading apostrophe so Excel keeps the exact text
# (otherwise it would silently coerce them to a number and drop things like
# trailing zeros / change "0.0720" -> 0.072). Prices that already contain two
# dots (e.g. "64.112.99") can never parse as a number, so no apostrophe is
# needed there. Percent cells (column E) are padded with literal spaces and
# a trailing "%", which Excel never autoconverts, so they are set as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.112.99"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "3.407.22"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'571.15"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").Value = "'162.75"
$ws.Range("E6").Value = "  +2.55%  "
$ws.Range("D8").Value = "3.407.85"
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("D9").Value = "'0.548"
$ws.Range("E9").Value = "  -4.20%  "
$ws.Range("E10").Value = "  +1.61%  "
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("D12").Value = "'0.419"
$ws.Range("E12").Value = "  -3.93%  "
$ws.Range("D13").Value = "3.994.36"
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").Value = "'26.86"
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").Value = "64.130.79"
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").Value = "3.401.60"
$ws.Range("E18").Value = "  -1.79%  "
$ws.Range("D19").Value = "'6.11"
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("D20").Value = "'13.48"
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").Value = "'372.01"
$ws.Range("E21").Value = "  -1.45%  "
$ws.Range("D22").Value = "'7.78"
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'70.17"
$ws.Range("E24").Value = "  -2.53%  "
$ws.Range("E25").Value = "  -3.37%  "
$ws.Range("D26").Value = "'0.0000115"
$ws.Range("E26").Value = "  -3.80%  "
$ws.Range("D27").Value = "'9.50"
$ws.Range("E27").Value = "  -3.83%  "
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "'6.09"
$ws.Range("E30").Value = "  +1.12%  "
$ws.Range("E31").Value = "  -2.59%  "
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("D34").Value = "'22.75"
$ws.Range("E34").Value = "  -1.61%  "
$ws.Range("D35").Value = "'7.01"
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("D36").Value = "'1.48"
$ws.Range("E36").Value = "  -5.23%  "
$ws.Range("D37").Value = "'159.52"
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("D38").Value = "'0.855"
$ws.Range("E38").Value = "  +8.59%  "
$ws.Range("E39").Value = "  -3.40%  "
$ws.Range("D40").Value = "'25.82"
$ws.Range("E40").Value = "  -0.51%  "
$ws.Range("D41").Value = "'0.0720"
$ws.Range("E41").Value = "  -3.06%  "
$ws.Range("D42").Value = "'42.72"
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'6.44"
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.731.10"
$ws.Range("E44").Value = "  -5.26%  "
$ws.Range("D45").Value = "'25.95"
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("D46").Value = "'4.35"
$ws.Range("E46").Value = "  -3.10%  "
$ws.Range("D47").Value = "'0.0304"
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("D49").Value = "'327.10"
$ws.Range("E49").Value = "  +2.16%  "
$ws.Range("E50").Value = "  -3.10%  "
$ws.Range("E51").Value = "  -1.65%  "
